# Apply the commit's changes:
#   1. Shared-string header text normalisation (same 3 strings on every
#      sheet's header row 1: A1, G1, I1).
#   2. On the "POSTRES" sheet, split the single B:H column-width band into
#      B:F / G / H so column G (the new text column) gets a wider custom
#      width.

$wb = $excel.ActiveWorkbook

# --- 1. Header text fixes, applied identically to every worksheet -------
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1").Value = "nombre"
    $ws.Range("G1").Value = "hidratos_de_carbono"
    $ws.Range("I1").Value = "composicion_en_equivalente"
}

# --- 2. Widen column G on the POSTRES sheet ------------------------------
$postres = $wb.Worksheets.Item("POSTRES")
$postres.Columns.Item(7).ColumnWidth = 27.43
